$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(15, 8).Value = 582.0172
$ws.Cells.Item(15, 9).Value = 582.0172
$ws.Cells.Item(15, 11).Value = 1746.0516
$ws.Cells.Item(15, 13).Value = -1577.0516

# ALC row 62
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(62, 8).Value = 6322.154
$ws.Cells.Item(62, 9).Value = 4480
$ws.Cells.Item(62, 10).Value = 6874.8
$ws.Cells.Item(62, 11).Value = 4480
$ws.Cells.Item(62, 12).Value = 6874.8
$ws.Cells.Item(62, 13).Value = -3856
$ws.Cells.Item(62, 14).Value = -8122.8

# ALC row 65
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(65, 8).Value = 6322.154
$ws.Cells.Item(65, 9).Value = 4480
$ws.Cells.Item(65, 10).Value = 6874.8
$ws.Cells.Item(65, 11).Value = 22400
$ws.Cells.Item(65, 12).Value = 34374
$ws.Cells.Item(65, 13).Value = -19280
$ws.Cells.Item(65, 14).Value = -40614

# ALC row 70
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(70, 8).Value = 89189.75
$ws.Cells.Item(70, 10).Value = 141203.8
$ws.Cells.Item(70, 12).Value = 423611.4
$ws.Cells.Item(70, 14).Value = -424151.4

# ALC row 73
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(73, 8).Value = 89189.75
$ws.Cells.Item(73, 10).Value = 141203.8
$ws.Cells.Item(73, 12).Value = 423611.4
$ws.Cells.Item(73, 14).Value = -425483.4

# ALC row 98
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(98, 8).Value = 1318.8182
$ws.Cells.Item(98, 9).Value = 1311.3334
$ws.Cells.Item(98, 11).Value = 1311.3334
$ws.Cells.Item(98, 13).Value = 186.6666

# ALC row 112
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(112, 8).Value = 2353.6428
$ws.Cells.Item(112, 10).Value = 2461.3333
$ws.Cells.Item(112, 12).Value = 7383.999899999999
$ws.Cells.Item(112, 14).Value = -9599.999899999999

# ALC row 122
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(122, 8).Value = 1318.8182
$ws.Cells.Item(122, 9).Value = 1311.3334
$ws.Cells.Item(122, 11).Value = 3934.0002
$ws.Cells.Item(122, 13).Value = -1484.0002

# ALC row 132
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(132, 8).Value = 2463.8572
$ws.Cells.Item(132, 9).Value = 2452.8696
$ws.Cells.Item(132, 11).Value = 7358.6088
$ws.Cells.Item(132, 13).Value = -4828.6088

# ALC row 137
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(137, 8).Value = 7732.227
$ws.Cells.Item(137, 10).Value = 13850.546
$ws.Cells.Item(137, 12).Value = 41551.638
$ws.Cells.Item(137, 14).Value = -46651.638

# ALC row 138
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(138, 8).Value = 5116.6
$ws.Cells.Item(138, 9).Value = 6157.9565
$ws.Cells.Item(138, 10).Value = 4607
$ws.Cells.Item(138, 11).Value = 18473.8695
$ws.Cells.Item(138, 12).Value = 13821
$ws.Cells.Item(138, 13).Value = -13333.8695
$ws.Cells.Item(138, 14).Value = -24101

# ARM row 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 1700.5714
$ws.Cells.Item(2, 9).Value = 1865.6364
$ws.Cells.Item(2, 10).Value = 1095.3334
$ws.Cells.Item(2, 11).Value = 1865.6364
$ws.Cells.Item(2, 12).Value = 1095.3334
$ws.Cells.Item(2, 13).Value = -1752.6364
$ws.Cells.Item(2, 14).Value = -1321.3334

# ARM row 32
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 18129.354
$ws.Cells.Item(32, 9).Value = 8018.7144
$ws.Cells.Item(32, 11).Value = 8018.7144
$ws.Cells.Item(32, 13).Value = -7731.7144

# ARM row 61
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 2398.8
$ws.Cells.Item(61, 9).Value = 2398
$ws.Cells.Item(61, 11).Value = 2398
$ws.Cells.Item(61, 13).Value = -2186

# ARM row 74
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value = 3723.08
$ws.Cells.Item(74, 9).Value = 1443.5333
$ws.Cells.Item(74, 11).Value = 1443.5333
$ws.Cells.Item(74, 13).Value = -569.5333000000001

# ARM row 77
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(77, 8).Value = 3723.08
$ws.Cells.Item(77, 9).Value = 1443.5333
$ws.Cells.Item(77, 11).Value = 7217.6665
$ws.Cells.Item(77, 13).Value = -2849.6665

# ARM row 116
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(116, 8).Value = 1700.5714
$ws.Cells.Item(116, 9).Value = 1865.6364
$ws.Cells.Item(116, 10).Value = 1095.3334
$ws.Cells.Item(116, 11).Value = 1865.6364
$ws.Cells.Item(116, 12).Value = 1095.3334
$ws.Cells.Item(116, 13).Value = 428.3635999999999
$ws.Cells.Item(116, 14).Value = -5683.3334

# ARM row 122
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(122, 8).Value = 306130.75
$ws.Cells.Item(122, 9).Value = 457105.38
$ws.Cells.Item(122, 11).Value = 1371316.14
$ws.Cells.Item(122, 13).Value = -1368866.14

# ARM row 132
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(132, 8).Value = 6255.567
$ws.Cells.Item(132, 9).Value = 1975.3572
$ws.Cells.Item(132, 11).Value = 5926.071599999999
$ws.Cells.Item(132, 13).Value = -3396.071599999999

# ARM row 136
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value = 2398.8
$ws.Cells.Item(136, 9).Value = 2398
$ws.Cells.Item(136, 11).Value = 7194
$ws.Cells.Item(136, 13).Value = -4644

# BSM row 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 1700.5714
$ws.Cells.Item(3, 9).Value = 1865.6364
$ws.Cells.Item(3, 10).Value = 1095.3334
$ws.Cells.Item(3, 11).Value = 1865.6364
$ws.Cells.Item(3, 12).Value = 1095.3334
$ws.Cells.Item(3, 13).Value = -1751.6364
$ws.Cells.Item(3, 14).Value = -1323.3334

# BSM row 64
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(64, 8).Value = 1471.5454
$ws.Cells.Item(64, 10).Value = 1698.8334
$ws.Cells.Item(64, 12).Value = 1698.8334
$ws.Cells.Item(64, 14).Value = -2148.8334

# BSM row 67
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(67, 8).Value = 1471.5454
$ws.Cells.Item(67, 10).Value = 1698.8334
$ws.Cells.Item(67, 12).Value = 1698.8334
$ws.Cells.Item(67, 14).Value = -3258.8334

# BSM row 105
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(105, 8).Value = 4021.2727
$ws.Cells.Item(105, 9).Value = 3224.682
$ws.Cells.Item(105, 10).Value = 5614.4546
$ws.Cells.Item(105, 11).Value = 3224.682
$ws.Cells.Item(105, 12).Value = 5614.4546
$ws.Cells.Item(105, 13).Value = -1477.682
$ws.Cells.Item(105, 14).Value = -9108.454600000001

# BSM row 134
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 5130.6665
$ws.Cells.Item(134, 9).Value = 3758.875
$ws.Cells.Item(134, 10).Value = 7874.25
$ws.Cells.Item(134, 11).Value = 11276.625
$ws.Cells.Item(134, 12).Value = 23622.75
$ws.Cells.Item(134, 13).Value = -8741.625
$ws.Cells.Item(134, 14).Value = -28692.75

# CRP row 134
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value = 2667.7727
$ws.Cells.Item(134, 9).Value = 1997.6316
$ws.Cells.Item(134, 11).Value = 5992.8948
$ws.Cells.Item(134, 13).Value = -3457.8948

# CUL row 2
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 8).Value = 254.33333
$ws.Cells.Item(2, 9).Value = 50
$ws.Cells.Item(2, 11).Value = 300
$ws.Cells.Item(2, 13).Value = -187

# CUL row 11
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(11, 8).Value = 391.64706
$ws.Cells.Item(11, 10).Value = 50
$ws.Cells.Item(11, 12).Value = 150
$ws.Cells.Item(11, 14).Value = -430

# CUL row 59
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(59, 8).Value = 8166
$ws.Cells.Item(59, 9).Value = 7454
$ws.Cells.Item(59, 11).Value = 22362
$ws.Cells.Item(59, 13).Value = -21822

# CUL row 107
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(107, 8).Value = 2573.375
$ws.Cells.Item(107, 10).Value = 1817.8
$ws.Cells.Item(107, 12).Value = 5453.4
$ws.Cells.Item(107, 14).Value = -9293.4

# GSM row 93
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()

# GSM row 102
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 2912.6667
$ws.Cells.Item(102, 9).Value = 1869
$ws.Cells.Item(102, 11).Value = 1869
$ws.Cells.Item(102, 13).Value = -247

# GSM row 113
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(113, 8).Value = 5000.909

# GSM row 122
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(122, 8).Value = 503930.53
$ws.Cells.Item(122, 9).Value = 66030.19
$ws.Cells.Item(122, 11).Value = 198090.57
$ws.Cells.Item(122, 13).Value = -195640.57

# GSM row 132
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value = 3140.5
$ws.Cells.Item(132, 9).Value = 2319.8462
$ws.Cells.Item(132, 11).Value = 6959.5386
$ws.Cells.Item(132, 13).Value = -4429.5386

# LTW row 22
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = 3019.4
$ws.Cells.Item(22, 9).Value = 1249
$ws.Cells.Item(22, 10).Value = 4199.6665
$ws.Cells.Item(22, 11).Value = 1249
$ws.Cells.Item(22, 12).Value = 4199.6665
$ws.Cells.Item(22, 13).Value = -954
$ws.Cells.Item(22, 14).Value = -4789.6665

# LTW row 27
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(27, 8).Value = 3019.4
$ws.Cells.Item(27, 9).Value = 1249
$ws.Cells.Item(27, 10).Value = 4199.6665
$ws.Cells.Item(27, 11).Value = 1249
$ws.Cells.Item(27, 12).Value = 4199.6665
$ws.Cells.Item(27, 13).Value = -1142
$ws.Cells.Item(27, 14).Value = -4413.6665

# LTW row 61
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(61, 8).Value = 3915.6924
$ws.Cells.Item(61, 9).Value = 3627.8635
$ws.Cells.Item(61, 10).Value = 5498.75
$ws.Cells.Item(61, 11).Value = 3627.8635
$ws.Cells.Item(61, 12).Value = 5498.75
$ws.Cells.Item(61, 13).Value = -3425.8635
$ws.Cells.Item(61, 14).Value = -5902.75

# LTW row 74
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()

# LTW row 77
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()

# LTW row 113
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(113, 8).Value = 3915.6924
$ws.Cells.Item(113, 9).Value = 3627.8635
$ws.Cells.Item(113, 10).Value = 5498.75
$ws.Cells.Item(113, 11).Value = 3627.8635
$ws.Cells.Item(113, 12).Value = 5498.75
$ws.Cells.Item(113, 13).Value = -1457.8635
$ws.Cells.Item(113, 14).Value = -9838.75

# LTW row 122
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(122, 8).Value = 7520.231
$ws.Cells.Item(122, 9).Value = 6382.1665
$ws.Cells.Item(122, 11).Value = 19146.4995
$ws.Cells.Item(122, 13).Value = -16696.4995

# LTW row 132
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(132, 8).Value = 5221.7896
$ws.Cells.Item(132, 9).Value = 3465.625
$ws.Cells.Item(132, 11).Value = 10396.875
$ws.Cells.Item(132, 13).Value = -7866.875

# LTW row 136
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(136, 8).Value = 3380.2
$ws.Cells.Item(136, 9).Value = 3380.2
$ws.Cells.Item(136, 11).Value = 10140.6
$ws.Cells.Item(136, 13).Value = -7590.599999999999

# WVR row 58
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(58, 8).Value = 9999.5
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 9999.5
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 12).Value = 9999.5
$ws.Cells.Item(58, 13).ClearContents()
$ws.Cells.Item(58, 14).Value = -10615.5

# WVR row 61
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(61, 8).Value = 8305
$ws.Cells.Item(61, 9).Value = 5000
$ws.Cells.Item(61, 10).Value = 9957.5
$ws.Cells.Item(61, 11).Value = 5000
$ws.Cells.Item(61, 12).Value = 9957.5
$ws.Cells.Item(61, 13).Value = -4708
$ws.Cells.Item(61, 14).Value = -10541.5

# WVR row 100
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(100, 8).Value = 2049.6667
$ws.Cells.Item(100, 10).Value = 1499.3334
$ws.Cells.Item(100, 12).Value = 2998.6668
$ws.Cells.Item(100, 14).Value = -4080.6668

# WVR row 122
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(122, 8).Value = 2200
$ws.Cells.Item(122, 9).Value = 2000.25
$ws.Cells.Item(122, 10).Value = 2999
$ws.Cells.Item(122, 11).Value = 6000.75
$ws.Cells.Item(122, 12).Value = 8997
$ws.Cells.Item(122, 13).Value = -3550.75
$ws.Cells.Item(122, 14).Value = -13897

# WVR row 132
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 2333.0527
$ws.Cells.Item(132, 9).Value = 1738
$ws.Cells.Item(132, 11).Value = 5214
$ws.Cells.Item(132, 13).Value = -2684
